# Update Katalon test-case execution timestamps (column B) to reflect the
# latest test run times, as captured in the commit
# "Added ABP test cases and modified IWP Bootstrap deferred test cases".

$wb = $excel.ActiveWorkbook

# PayNowCC: rows 2-7 hold the per test-case "last executed" timestamp.
$ws = $wb.Worksheets.Item("PayNowCC")
$ws.Range("B2").Value = "Tue Sep 02 11:33:59 IST 2025"
$ws.Range("B3").Value = "Tue Sep 02 11:34:49 IST 2025"
$ws.Range("B4").Value = "Tue Sep 02 11:35:36 IST 2025"
$ws.Range("B5").Value = "Tue Sep 02 11:36:25 IST 2025"
$ws.Range("B6").Value = "Tue Sep 02 11:37:16 IST 2025"
$ws.Range("B7").Value = "Tue Sep 02 11:38:03 IST 2025"

$ws = $wb.Worksheets.Item("CCDeferredCC")
$ws.Range("B2").Value = "Tue Sep 02 11:09:39 IST 2025"

$ws = $wb.Worksheets.Item("CMCAutopayCC")
$ws.Range("B2").Value = "Tue Sep 02 11:15:27 IST 2025"

$ws = $wb.Worksheets.Item("NoModifyAmountCC")
$ws.Range("B2").Value = "Tue Sep 02 11:22:05 IST 2025"

$ws = $wb.Worksheets.Item("NoModifyBillingAddressCC")
$ws.Range("B2").Value = "Tue Sep 02 11:26:20 IST 2025"

$ws = $wb.Worksheets.Item("PayNowCreditCardDCF")
$ws.Range("B2").Value = "Tue Sep 02 11:30:22 IST 2025"

$ws = $wb.Worksheets.Item("PayNowCreditCardSCF")
$ws.Range("B2").Value = "Tue Sep 02 11:43:03 IST 2025"

$ws = $wb.Worksheets.Item("DCFCCVerbiage")
$ws.Range("B2").Value = "Tue Sep 02 11:46:31 IST 2025"

$ws = $wb.Worksheets.Item("SCFCCVerbiage")
$ws.Range("B2").Value = "Tue Sep 02 11:48:44 IST 2025"
